$d = $word.ActiveDocument

# Locate the target paragraph ("Once we have the String command...").
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Once we have the String command, we enter into a switch-case*") {
        $target = $p
        break
    }
}

# The paragraph ends with a hidden "_GoBack" bookmark sitting right at the
# end of the run text (before the paragraph mark). Remove it first so the
# upcoming InsertXML (which replaces the whole range) doesn't leave it
# stranded at the front of the new content; we re-emit it ourselves at the
# tail end of the replacement XML, exactly where it originally sat.
$bookmark = $d.Bookmarks("_GoBack")
$bookmark.Delete()

# Range covering the paragraph's run text, excluding the trailing
# paragraph mark.
$r = $d.Range($target.Range.Start, $target.Range.End - 1)

$runsXml = '<w:r><w:t xml:space="preserve">After that, we have two big nested switch-cases. Our logic here is the following. If the user has pressed one of the sidebar menu buttons, then all the application would have to do is reload the main area of the dashboard with a new page. Hence, in this </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>case</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> we set the command to be &#8220;</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>reloadMain</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>&#8221; and we also add another attribute to the request, called &#8220;</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>mainArea</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>&#8221;</w:t></w:r><w:r><w:t xml:space="preserve">, which holds information on which </w:t></w:r><w:r><w:t>.</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>jsp</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> should be shown</w:t></w:r><w:r><w:t xml:space="preserve"> in the main area</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r><w:t xml:space="preserve"> Anything else the user wants to do is a different command and </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>is handled</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> by corresponding methods.</w:t></w:r>'

$bookmarkXml = '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $runsXml + $bookmarkXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$r.InsertXML($xml)
